$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value2 = $ws.Range("F14").Value2
$ws.Range("F14").Clear() | Out-Null
$ws.Range("F13").Clear() | Out-Null

$ws.Range("D23").Clear() | Out-Null
$ws.Range("D22").Clear() | Out-Null

$ws.Columns("C").ColumnWidth = 17.0
$ws.Columns("D").ColumnWidth = 23.666666666666668
$ws.Columns("E").ColumnWidth = 20.666666666666668

$ws.Range("E20").Select() | Out-Null
